$d = $word.ActiveDocument

# 1) Remove the "_GoBack" bookmark (bookmarkStart/bookmarkEnd pair)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Add a tab run after the "Задание 02" run, i.e. right after the bold/underlined
#    run containing "2" following "Задание 0".
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Задание 02", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPoint = $d.Range($rng.End, $rng.End)
    $insertPoint.InsertAfter("`t")
}
